$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.060.23'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.67%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.404.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.78%  '

# Row 4
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.75%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '665.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.89%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.44'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.69%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.431'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.85%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.05'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.33%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.05%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.401.57'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.79%  '

# Row 12
$ws.Range("E12").Value = '  +3.14%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.40'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.71%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +15.91%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '97.762.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.27%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000265'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.32%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.040.51'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.75%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.04'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +21.85%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.410.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.05%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.558'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +28.76%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.84%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.89%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.43'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.03%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '508.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.51%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000205'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.13%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.88%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '100.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.76%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.27%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.597.85'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.27%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.148'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.26%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.57'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.62%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.198'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.15%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.17%  '

# Row 34
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.02%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +15.47%  '

# Row 36
$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.572'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.30%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '29.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.03%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.51'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.63%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.47%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '534.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.66%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.152'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.02%  '

# Row 42
$ws.Range("E42").Value = '  +0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.875'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.95%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.14%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +17.42%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +18.34%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0426'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.34%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +16.63%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.91%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.34%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.63'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.93%  '
